$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new schema row, "discount_upto", is introduced right above the
# existing "item_constraint" row (row 8). Inserting an entire row there
# pushes "item_constraint" and everything below it (minimum_amount,
# customer_id, campaign_external_id, visibility, coupon_id, campaign_id,
# expires_at) down by one row, carrying their values/styles with them --
# exactly what a real "insert row" edit in Excel does.
$ws.Rows.Item(8).Insert()

# Fill in the new row 8 with the "discount_upto" field definition:
# key | null_allowed | type | ... | examples
$ws.Range("A8").Value = "discount_upto"
$ws.Range("B8").Value = "Yes"
$ws.Range("C8").Value = "number"
$ws.Range("K8").Value = 100.0
